$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: advance the date shown in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value2 = $ws.Range("A1").Value2 + 1

# Step 2: update the price list for FRATACHO de Pino items (rows 33-36)
$ws.Range("D33").Value2 = 344
$ws.Range("D34").Value2 = 393
$ws.Range("D35").Value2 = 424
$ws.Range("D36").Value2 = 465
